# Update NATMI Dkk1-Lrp5 TPM-derived numbers (commit: "update scripts wuth new tpm")
#
# The underlying receptor (Lrp5) expression data used to compute this
# ligand-receptor pair table was refreshed with new TPM values. That changes
# the "Receptor total expression value" feeding row 2 (target cluster ECs),
# which ripples into every column derived from it (average expression,
# edge weights, and the specificity columns that are normalised across the
# three target-cluster rows for this Dkk1/Lrp5 pair). Rows 3 and 4 keep the
# same receptor totals, but their specificity columns still shift because the
# normalisation denominator (summed across rows 2-4) changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Sending cluster MuSCs, Dkk1 -> Lrp5, Target cluster ECs) ---
$ws.Range("H2").Value = 0.09138
$ws.Range("M2").Value = 15.01856033333333
$ws.Range("N2").Value = 45.055681
$ws.Range("O2").Value = 0.4908713633047416
$ws.Range("P2").Value = 0.4908713633047417
$ws.Range("Q2").Value = 0.4574653477533334
$ws.Range("R2").Value = 4.11718812978
$ws.Range("S2").Value = 0.4908713633047416
$ws.Range("T2").Value = 0.4908713633047417

# --- Row 3 (Target cluster FAPs) ---
$ws.Range("H3").Value = 0.09138
$ws.Range("O3").Value = 0.3099803572711625
$ws.Range("P3").Value = 0.3099803572711625
$ws.Range("Q3").Value = 0.2888847925066667
$ws.Range("S3").Value = 0.3099803572711625
$ws.Range("T3").Value = 0.3099803572711625

# --- Row 4 (Target cluster MuSCs) ---
$ws.Range("H4").Value = 0.09138
$ws.Range("O4").Value = 0.1991482794240958
$ws.Range("P4").Value = 0.1991482794240958
$ws.Range("S4").Value = 0.1991482794240958
$ws.Range("T4").Value = 0.1991482794240958
